$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.197692513465881
$ws.Range("B1").Value = 2.31177830696106
$ws.Range("C1").Value = 3.381871461868286
$ws.Range("D1").Value = 3.397888660430908
$ws.Range("E1").Value = 1.140867352485657
